$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-10-24 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-10-25 Friday", 2) | Out-Null
$d.Content.Find.Execute("862÷4=215, 2", $true, $false, $false, $false, $false, $true, 1, $false, "205÷7=29, 2", 2) | Out-Null
$d.Content.Find.Execute("510÷6=85, 0", $true, $false, $false, $false, $false, $true, 1, $false, "746÷2=373, 0", 2) | Out-Null
$d.Content.Find.Execute("868÷3=289, 1", $true, $false, $false, $false, $false, $true, 1, $false, "112÷7=16, 0", 2) | Out-Null
$d.Content.Find.Execute("118÷8=14, 6", $true, $false, $false, $false, $false, $true, 1, $false, "548÷3=182, 2", 2) | Out-Null
$d.Content.Find.Execute("318÷6=53, 0", $true, $false, $false, $false, $false, $true, 1, $false, "445÷5=89, 0", 2) | Out-Null
$d.Content.Find.Execute("700÷3=233, 1", $true, $false, $false, $false, $false, $true, 1, $false, "182÷8=22, 6", 2) | Out-Null
$d.Content.Find.Execute("789÷7=112, 5", $true, $false, $false, $false, $false, $true, 1, $false, "837÷5=167, 2", 2) | Out-Null
$d.Content.Find.Execute("944÷7=134, 6", $true, $false, $false, $false, $false, $true, 1, $false, "945÷3=315, 0", 2) | Out-Null
$d.Content.Find.Execute("889÷7=127, 0", $true, $false, $false, $false, $false, $true, 1, $false, "466÷2=233, 0", 2) | Out-Null
$d.Content.Find.Execute("957÷8=119, 5", $true, $false, $false, $false, $false, $true, 1, $false, "667÷6=111, 1", 2) | Out-Null
$d.Content.Find.Execute("320÷8=40, 0", $true, $false, $false, $false, $false, $true, 1, $false, "908÷7=129, 5", 2) | Out-Null
$d.Content.Find.Execute("609÷7=87, 0", $true, $false, $false, $false, $false, $true, 1, $false, "528÷4=132, 0", 2) | Out-Null
$d.Content.Find.Execute("258÷4=64, 2", $true, $false, $false, $false, $false, $true, 1, $false, "864÷3=288, 0", 2) | Out-Null
$d.Content.Find.Execute("131÷8=16, 3", $true, $false, $false, $false, $false, $true, 1, $false, "620÷6=103, 2", 2) | Out-Null
$d.Content.Find.Execute("878÷4=219, 2", $true, $false, $false, $false, $false, $true, 1, $false, "499÷8=62, 3", 2) | Out-Null
$d.Content.Find.Execute("313÷5=62, 3", $true, $false, $false, $false, $false, $true, 1, $false, "330÷8=41, 2", 2) | Out-Null
$d.Content.Find.Execute("701÷6=116, 5", $true, $false, $false, $false, $false, $true, 1, $false, "436÷3=145, 1", 2) | Out-Null
$d.Content.Find.Execute("298÷2=149, 0", $true, $false, $false, $false, $false, $true, 1, $false, "705÷2=352, 1", 2) | Out-Null
$d.Content.Find.Execute("106÷2=53, 0", $true, $false, $false, $false, $false, $true, 1, $false, "268÷9=29, 7", 2) | Out-Null
$d.Content.Find.Execute("651÷7=93, 0", $true, $false, $false, $false, $false, $true, 1, $false, "224÷6=37, 2", 2) | Out-Null
$d.Content.Find.Execute("999÷2=499, 1", $true, $false, $false, $false, $false, $true, 1, $false, "952÷4=238, 0", 2) | Out-Null
$d.Content.Find.Execute("354÷3=118, 0", $true, $false, $false, $false, $false, $true, 1, $false, "847÷4=211, 3", 2) | Out-Null
$d.Content.Find.Execute("271÷5=54, 1", $true, $false, $false, $false, $false, $true, 1, $false, "303÷3=101, 0", 2) | Out-Null
$d.Content.Find.Execute("290÷7=41, 3", $true, $false, $false, $false, $false, $true, 1, $false, "806÷7=115, 1", 2) | Out-Null
$d.Content.Find.Execute("511÷9=56, 7", $true, $false, $false, $false, $false, $true, 1, $false, "319÷9=35, 4", 2) | Out-Null

Write-Output "done"
